$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing "*" from the condition labels (column A, rows 2-6)
$ws.Range("A2").Value = "That recipient countries comply with`nclimate targets and human rights"
$ws.Range("A3").Value = "That recipient countries cooperate`nto fight illegal migrations"
$ws.Range("A4").Value = "That other high-income countries`nalso increase their foreign aid"
$ws.Range("A5").Value = "That this is financed by increased taxes on millionaires"
$ws.Range("A6").Value = "That we can be sure the aid reaches`npeople in need and money is not diverted"

# Re-assigning multi-line text can trigger Excel's row auto-height; restore
# the rows to their original (default) auto-fit height so no explicit
# row height ends up stored in the worksheet.
$ws.Rows("2").RowHeight = 15
$ws.Rows("2").AutoFit()
$ws.Rows("3").RowHeight = 15
$ws.Rows("3").AutoFit()
$ws.Rows("4").RowHeight = 15
$ws.Rows("4").AutoFit()
$ws.Rows("6").RowHeight = 15
$ws.Rows("6").AutoFit()

# Overwrite the figures in column B with corrected precision values
$ws.Range("B2").Value = 0.612607348826151
$ws.Range("B3").Value = 0.355027957335685
$ws.Range("B4").Value = 0.446393661326311
$ws.Range("B5").Value = 0.363253941831613
$ws.Range("B6").Value = 0.684337117444274
